$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename/retitle existing track entries (credits/assets edits)
$ws.Range("B21").Value = "Renegade Level 2"
$ws.Range("B23").Value = "Renegade Level 2"
$ws.Range("B27").Value = "Escape (Runa)"
$ws.Range("B20").Value = "Anima Chant (Runa)"

# New track added for Ch26 (Runa)
$ws.Range("B37").Value = "LoZ Twilight Princess Midna's Lament"

# Update the active selection to match the latest edit location
$ws.Range("B38").Select() | Out-Null
